$d = $word.ActiveDocument

# Update the date paragraph at the top of the document
$dateRange = $d.Paragraphs.Item(1).Range
$dateRange.Find.Execute("2025-12-11 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-12 Friday", 2) | Out-Null

# Update the division problems in the table, cell by cell
$tbl = $d.Tables.Item(1)

$cellRange = $tbl.Cell(1, 1).Range
$cellRange.Find.Execute("12÷8=1, 4", $true, $false, $false, $false, $false, $true, 1, $false, "66÷4=16, 2", 2) | Out-Null
$cellRange = $tbl.Cell(1, 2).Range
$cellRange.Find.Execute("56÷5=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "32÷7=4, 4", 2) | Out-Null
$cellRange = $tbl.Cell(1, 3).Range
$cellRange.Find.Execute("41÷9=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "22÷8=2, 6", 2) | Out-Null
$cellRange = $tbl.Cell(1, 4).Range
$cellRange.Find.Execute("62÷3=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "15÷7=2, 1", 2) | Out-Null
$cellRange = $tbl.Cell(1, 5).Range
$cellRange.Find.Execute("81÷9=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "11÷9=1, 2", 2) | Out-Null
$cellRange = $tbl.Cell(5, 1).Range
$cellRange.Find.Execute("43÷9=4, 7", $true, $false, $false, $false, $false, $true, 1, $false, "88÷5=17, 3", 2) | Out-Null
$cellRange = $tbl.Cell(5, 2).Range
$cellRange.Find.Execute("44÷8=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "90÷8=11, 2", 2) | Out-Null
$cellRange = $tbl.Cell(5, 3).Range
$cellRange.Find.Execute("38÷6=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "55÷7=7, 6", 2) | Out-Null
$cellRange = $tbl.Cell(5, 4).Range
$cellRange.Find.Execute("25÷8=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "71÷2=35, 1", 2) | Out-Null
$cellRange = $tbl.Cell(5, 5).Range
$cellRange.Find.Execute("21÷7=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "16÷6=2, 4", 2) | Out-Null
$cellRange = $tbl.Cell(9, 1).Range
$cellRange.Find.Execute("35÷3=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "83÷5=16, 3", 2) | Out-Null
$cellRange = $tbl.Cell(9, 2).Range
$cellRange.Find.Execute("20÷5=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "72÷4=18, 0", 2) | Out-Null
$cellRange = $tbl.Cell(9, 3).Range
$cellRange.Find.Execute("84÷7=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "59÷7=8, 3", 2) | Out-Null
$cellRange = $tbl.Cell(9, 4).Range
$cellRange.Find.Execute("33÷2=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "25÷5=5, 0", 2) | Out-Null
$cellRange = $tbl.Cell(9, 5).Range
$cellRange.Find.Execute("90÷4=22, 2", $true, $false, $false, $false, $false, $true, 1, $false, "33÷3=11, 0", 2) | Out-Null
$cellRange = $tbl.Cell(13, 1).Range
$cellRange.Find.Execute("22÷7=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "80÷2=40, 0", 2) | Out-Null
$cellRange = $tbl.Cell(13, 2).Range
$cellRange.Find.Execute("29÷6=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "75÷4=18, 3", 2) | Out-Null
$cellRange = $tbl.Cell(13, 3).Range
$cellRange.Find.Execute("14÷7=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "27÷7=3, 6", 2) | Out-Null
$cellRange = $tbl.Cell(13, 4).Range
$cellRange.Find.Execute("43÷7=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "76÷4=19, 0", 2) | Out-Null
$cellRange = $tbl.Cell(13, 5).Range
$cellRange.Find.Execute("63÷7=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "57÷5=11, 2", 2) | Out-Null
$cellRange = $tbl.Cell(17, 1).Range
$cellRange.Find.Execute("65÷2=32, 1", $true, $false, $false, $false, $false, $true, 1, $false, "59÷6=9, 5", 2) | Out-Null
$cellRange = $tbl.Cell(17, 2).Range
$cellRange.Find.Execute("33÷5=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "64÷2=32, 0", 2) | Out-Null
$cellRange = $tbl.Cell(17, 3).Range
$cellRange.Find.Execute("13÷6=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "41÷4=10, 1", 2) | Out-Null
$cellRange = $tbl.Cell(17, 4).Range
$cellRange.Find.Execute("50÷8=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "35÷5=7, 0", 2) | Out-Null
$cellRange = $tbl.Cell(17, 5).Range
$cellRange.Find.Execute("52÷3=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "21÷7=3, 0", 2) | Out-Null
